$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.318.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3952"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07952"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.935.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.118"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.779"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06952"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.323.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.365"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.152.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.065"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.125"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.002"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09392"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9312"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.361"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.362"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.278"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.213"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05844"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02113"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.62%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5762"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.304"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5433"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07080"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.570"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
